# Auto-generated edit script: applies cell-value updates described by the commit diff
# to the "Behemoth_Profits" workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 502.33334
$ws.Range("J17").Value = 502.33334
$ws.Range("L17").Value = 1507.00002
$ws.Range("N17").Value = -1843.00002
$ws.Range("H40").Value = 3536.04
$ws.Range("I40").Value = 2938.5386
$ws.Range("K40").Value = 2938.5386
$ws.Range("M40").Value = -2763.5386
$ws.Range("H69").Value = 19436
$ws.Range("I69").Value = 9245.5
$ws.Range("K69").Value = 27736.5
$ws.Range("M69").Value = -26862.5
$ws.Range("H70").Value = 2204.3333
$ws.Range("J70").Value = 2665.8333
$ws.Range("L70").Value = 7997.499899999999
$ws.Range("N70").Value = -8537.499899999999
$ws.Range("H72").Value = 19436
$ws.Range("I72").Value = 9245.5
$ws.Range("K72").Value = 83209.5
$ws.Range("M72").Value = -78841.5
$ws.Range("H73").Value = 2204.3333
$ws.Range("J73").Value = 2665.8333
$ws.Range("L73").Value = 7997.499899999999
$ws.Range("N73").Value = -9869.499899999999
$ws.Range("H88").Value = 674773.25
$ws.Range("I88").Value = 1624.5
$ws.Range("K88").Value = 1624.5
$ws.Range("M88").Value = -1218.5
$ws.Range("H91").Value = 674773.25
$ws.Range("I91").Value = 1624.5
$ws.Range("K91").Value = 1624.5
$ws.Range("M91").Value = -220.5
$ws.Range("H107").Value = 1189.9445
$ws.Range("I107").Value = 1209.7693
$ws.Range("J107").Value = 1138.4
$ws.Range("K107").Value = 1209.7693
$ws.Range("L107").Value = 1138.4
$ws.Range("M107").Value = 710.2307000000001
$ws.Range("N107").Value = -4978.4
$ws.Range("H116").Value = 6599
$ws.Range("I116").Value = 6071.2856
$ws.Range("K116").Value = 6071.2856
$ws.Range("M116").Value = -2629.2856
$ws.Range("H137").Value = 5877
$ws.Range("I137").Value = 3588.75
$ws.Range("J137").Value = 11368.8
$ws.Range("K137").Value = 10766.25
$ws.Range("L137").Value = 34106.39999999999
$ws.Range("M137").Value = -8216.25
$ws.Range("N137").Value = -39206.39999999999
$ws.Range("H138").Value = 1236866.4
$ws.Range("I138").Value = 921.36365
$ws.Range("J138").Value = 2086578.6
$ws.Range("K138").Value = 2764.09095
$ws.Range("L138").Value = 6259735.800000001
$ws.Range("M138").Value = 2375.90905
$ws.Range("N138").Value = -6270015.800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 73.666664
$ws.Range("I5").Value = 73.666664
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 73.666664
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = 38.333336
$ws.Range("H61").Value = 45552620
$ws.Range("I61").Value = 125003670
$ws.Range("J61").Value = 152015.14
$ws.Range("K61").Value = 125003670
$ws.Range("L61").Value = 152015.14
$ws.Range("M61").Value = -125003458
$ws.Range("N61").Value = -152439.14
$ws.Range("H74").Value = 8937252
$ws.Range("I74").Value = 13159837
$ws.Range("J74").Value = 22906.666
$ws.Range("K74").Value = 13159837
$ws.Range("L74").Value = 22906.666
$ws.Range("M74").Value = -13158963
$ws.Range("N74").Value = -24654.666
$ws.Range("H77").Value = 8937252
$ws.Range("I77").Value = 13159837
$ws.Range("J77").Value = 22906.666
$ws.Range("K77").Value = 65799185
$ws.Range("L77").Value = 114533.33
$ws.Range("M77").Value = -65794817
$ws.Range("N77").Value = -123269.33
$ws.Range("H88").Value = 1926.3636
$ws.Range("I88").Value = 1898.75
$ws.Range("J88").Value = 1942.1428
$ws.Range("K88").Value = 1898.75
$ws.Range("L88").Value = 1942.1428
$ws.Range("M88").Value = -1492.75
$ws.Range("N88").Value = -2754.1428
$ws.Range("H91").Value = 1926.3636
$ws.Range("I91").Value = 1898.75
$ws.Range("J91").Value = 1942.1428
$ws.Range("K91").Value = 1898.75
$ws.Range("L91").Value = 1942.1428
$ws.Range("M91").Value = -494.75
$ws.Range("N91").Value = -4750.1428
$ws.Range("H132").Value = 6140.645
$ws.Range("I132").Value = 3757.5
$ws.Range("K132").Value = 11272.5
$ws.Range("M132").Value = -8742.5
$ws.Range("H136").Value = 45552620
$ws.Range("I136").Value = 125003670
$ws.Range("J136").Value = 152015.14
$ws.Range("K136").Value = 375011010
$ws.Range("L136").Value = 456045.42
$ws.Range("M136").Value = -375008460
$ws.Range("N136").Value = -461145.42

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 73.666664
$ws.Range("I4").Value = 73.666664
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 73.666664
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = 41.333336
$ws.Range("H86").Value = 2076.5625
$ws.Range("I86").Value = 2144.0833
$ws.Range("J86").Value = 1874
$ws.Range("K86").Value = 2144.0833
$ws.Range("L86").Value = 1874
$ws.Range("M86").Value = -1021.0833
$ws.Range("N86").Value = -4120
$ws.Range("H89").Value = 2076.5625
$ws.Range("I89").Value = 2144.0833
$ws.Range("J89").Value = 1874
$ws.Range("K89").Value = 10720.4165
$ws.Range("L89").Value = 9370
$ws.Range("M89").Value = -5104.416499999999
$ws.Range("N89").Value = -20602
$ws.Range("H99").Value = 2866.0625
$ws.Range("I99").Value = 2321
$ws.Range("K99").Value = 2321
$ws.Range("M99").Value = -823
$ws.Range("H107").Value = 1668.7
$ws.Range("I107").Value = 1396.75
$ws.Range("J107").Value = 2756.5
$ws.Range("K107").Value = 1396.75
$ws.Range("L107").Value = 2756.5
$ws.Range("M107").Value = 523.25
$ws.Range("N107").Value = -6596.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 999.5
$ws.Range("I16").Value = 999.5
$ws.Range("K16").Value = 999.5
$ws.Range("M16").Value = -712.5
$ws.Range("H31").Value = 779640.4399999999
$ws.Range("I31").Value = 1472.1333
$ws.Range("J31").Value = 1557808.8
$ws.Range("K31").Value = 1472.1333
$ws.Range("L31").Value = 1557808.8
$ws.Range("M31").Value = -1177.1333
$ws.Range("N31").Value = -1558398.8
$ws.Range("H34").Value = 779640.4399999999
$ws.Range("I34").Value = 1472.1333
$ws.Range("J34").Value = 1557808.8
$ws.Range("K34").Value = 1472.1333
$ws.Range("L34").Value = 1557808.8
$ws.Range("M34").Value = -1270.1333
$ws.Range("N34").Value = -1558212.8
$ws.Range("H39").Value = 5000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 5000
$ws.Range("K39").Value = 0
$ws.Range("L39").ClearContents()
$ws.Range("M39").Value = 5000
$ws.Range("N39").Value = -5782
$ws.Range("H49").Value = 5000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 5000
$ws.Range("K49").Value = 0
$ws.Range("L49").ClearContents()
$ws.Range("M49").Value = 5000
$ws.Range("N49").Value = -5364
$ws.Range("H107").Value = 797.0625
$ws.Range("I107").Value = 479
$ws.Range("J107").Value = 2175.3333
$ws.Range("K107").Value = 479
$ws.Range("L107").Value = 2175.3333
$ws.Range("M107").Value = 1441
$ws.Range("N107").Value = -6015.3333
$ws.Range("H113").Value = 999.5
$ws.Range("I113").Value = 999.5
$ws.Range("K113").Value = 999.5
$ws.Range("M113").Value = 1170.5
$ws.Range("H134").Value = 669279.75
$ws.Range("I134").Value = 770553.1
$ws.Range("K134").Value = 2311659.3
$ws.Range("M134").Value = -2309124.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8605328
$ws.Range("I4").Value = 8000044.5
$ws.Range("K4").Value = 24000133.5
$ws.Range("M4").Value = -24000021.5
$ws.Range("H129").Value = 27861538
$ws.Range("J129").Value = 37148416
$ws.Range("L129").Value = 111445248
$ws.Range("N129").Value = -111455248
$ws.Range("H141").Value = 12518.6
$ws.Range("J141").Value = 14488.889
$ws.Range("L141").Value = 43466.667
$ws.Range("N141").Value = -53826.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 16111.2
$ws.Range("I5").Value = 10137.75
$ws.Range("K5").Value = 10137.75
$ws.Range("M5").Value = -10025.75
$ws.Range("H52").Value = 50000
$ws.Range("I52").Value = 50000
$ws.Range("K52").Value = 50000
$ws.Range("M52").Value = -49741
$ws.Range("H57").Value = 12500.5
$ws.Range("I57").Value = 12500.5
$ws.Range("K57").Value = 12500.5
$ws.Range("M57").Value = -11680.5
$ws.Range("H132").Value = 66668730
$ws.Range("I132").Value = 111112530
$ws.Range("K132").Value = 333337590
$ws.Range("M132").Value = -333335060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2790.8
$ws.Range("I46").Value = 2553.5
$ws.Range("K46").Value = 2553.5
$ws.Range("M46").Value = -2365.5
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H100").Value = 6237.3335
$ws.Range("I100").Value = 6884.8
$ws.Range("K100").Value = 6884.8
$ws.Range("M100").Value = -6343.8
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 38495
$ws.Range("J54").Value = 38495
$ws.Range("L54").Value = 38495
$ws.Range("N54").Value = -39535
$ws.Range("H122").Value = 4482.448
$ws.Range("I122").Value = 2442.8096
$ws.Range("J122").Value = 9836.5
$ws.Range("K122").Value = 7328.4288
$ws.Range("L122").Value = 29509.5
$ws.Range("M122").Value = -4878.4288
$ws.Range("N122").Value = -34409.5
$ws.Range("H126").Value = 7191.4546
$ws.Range("I126").Value = 3786.5715
$ws.Range("K126").Value = 11359.7145
$ws.Range("M126").Value = -8889.7145
